# sm_car_data_Aero_Coefficients.xlsx — "Update 2p0. Convention change to
# support multi-axle vehicles"
#
# Adds two new vehicle-instance sheets:
#   - Truck_Amandla   (inserted right after Bus_Makhulu)
#   - Trailer_Kumanzi (appended at the end, becomes the active sheet)
#
# Both new sheets are built the same way the original author built them:
# duplicate the existing "Bus_Makhulu" sheet (same CD/CL/rho/ARef values,
# same "sedan" class marker) and then only change the Instance name (H3)
# and the sPressureCentre row (F9:H9).

$wb = $excel.ActiveWorkbook

# --- Truck_Amandla: copy of Bus_Makhulu, placed immediately after it ---
$bus = $wb.Worksheets.Item("Bus_Makhulu")
$bus.Copy($null, $bus)
$truck = $wb.Worksheets.Item($bus.Index + 1)
$truck.Name = "Truck_Amandla"
$truck.Range("H3").Value = "Truck_Amandla"
$truck.Range("F9").Value = -1.2
$truck.Range("G9").Value = 0
$truck.Range("H9").Value = 1.1
$truck.Activate()
[void]$truck.Range("H5:H9").Select()

# --- Trailer_Kumanzi: copy of Bus_Makhulu, appended as the last sheet ---
$bus2 = $wb.Worksheets.Item("Bus_Makhulu")
$bus2.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$kumanzi = $wb.Worksheets.Item($wb.Worksheets.Count)
$kumanzi.Name = "Trailer_Kumanzi"
$kumanzi.Range("H3").Value = "Trailer_Kumanzi"
$kumanzi.Range("F9").Value = 5
$kumanzi.Range("G9").Value = 0
$kumanzi.Range("H9").Value = 2

# Trailer_Kumanzi ends up as the active/selected sheet, matching the saved
# workbook state in the commit.
$kumanzi.Activate()
[void]$kumanzi.Range("J20").Select()
